# COSTA-RICA - RBM Impact Indicators from ENAHO
# Adding documentation and format:
#   - bump the "last updated" date on the title slide
#   - fix wording on the statistics-inclusion slide
#   - nudge a batch of auto-generated chart labels (DejaVu Sans -> Arial,
#     slightly tighter bounding boxes) on the three outcome/impact slides

$p = $ppt.ActivePresentation

# EMU -> points helper. The host truncates (rather than rounds) when it
# converts the Single-precision point value back to EMU, so nudge the
# value up by half an EMU before dividing to land on the exact integer.
function EmuToPt($emu) {
    return ([double]$emu + 0.5) / 12700.0
}

function SetShapeBox($container, $name, $x, $y, $cx, $cy) {
    $sh = $container.GroupItems.Item($name)
    $sh.Left   = EmuToPt $x
    $sh.Top    = EmuToPt $y
    $sh.Width  = EmuToPt $cx
    $sh.Height = EmuToPt $cy
}

function SetShapeFont($container, $name, $fontName) {
    $sh = $container.GroupItems.Item($name)
    $tr = $sh.TextFrame.TextRange
    $tr.Font.Name = $fontName
    $tr.Font.NameComplexScript = $fontName
}

# ---------------------------------------------------------------------
# Slide 1 - title slide date placeholder
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)
$dateShape = $s1.Shapes.Item(3)
$dateShape.TextFrame.TextRange.Text = "21 November 2022"

# ---------------------------------------------------------------------
# Slide 2 - "Inclusión estadística" body copy wording fixes
# ---------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$body = $s2.Shapes.Item(2)
$tr2 = $body.TextFrame.TextRange

# Apply the later edit first so the earlier (length-changing) edit does
# not shift the character offsets of the one that follows it.
$thirdParagraphRun = $tr2.Characters(498, 242)
$thirdParagraphRun.Text = "El ACNUR dispone de un conjunto de indicadores básicos. En Costa Rica, la Encuesta Nacional de Hogares incluye a las personas desplazadas por la fuerza y el ACNUR puede utilizar las estadísticas oficiales para calcular sus propios indicadores."

$firstParagraphTailRun = $tr2.Characters(97, 143)
$firstParagraphTailRun.Text = " de 2018,las Oficinas Nacionales de Estadística son las mejores posicionadas para producir datos de alta calidad sobre el desplazamiento forzado."

# ---------------------------------------------------------------------
# Slides 4-6 - auto-generated chart label textboxes (inside the "grp"
# group shape that is the 2nd top-level shape on each slide)
# ---------------------------------------------------------------------

# Slide 4 - Impact Area 2
$s4 = $p.Slides.Item(4)
$g4 = $s4.Shapes.Item(2)
SetShapeBox  $g4 "tx17" 5657183 3515857 397540 105768
SetShapeFont $g4 "tx17" "Arial"
SetShapeBox  $g4 "tx18" 6273249 3189685 397540 105768
SetShapeFont $g4 "tx18" "Arial"
SetShapeBox  $g4 "tx19" 10590864 4060893 397540 105768
SetShapeFont $g4 "tx19" "Arial"
SetShapeBox  $g4 "tx20" 10193115 3734720 397540 105768
SetShapeFont $g4 "tx20" "Arial"
SetShapeBox  $g4 "tx21" 10980922 4605928 397540 105768
SetShapeFont $g4 "tx21" "Arial"
SetShapeBox  $g4 "tx22" 8329753 4279755 397540 105768
SetShapeFont $g4 "tx22" "Arial"
SetShapeBox  $g4 "tx36" 3727579 2764413 77974 75989
SetShapeFont $g4 "tx36" "Arial"
SetShapeBox  $g4 "tx38" 5051978 2764413 77974 75989
SetShapeFont $g4 "tx38" "Arial"

# Slide 5 - Impact Area 3
$s5 = $p.Slides.Item(5)
$g5 = $s5.Shapes.Item(2)
SetShapeBox  $g5 "tx14" 10628300 3664503 397540 105768
SetShapeFont $g5 "tx14" "Arial"
SetShapeBox  $g5 "tx15" 9897464 3338331 397540 105768
SetShapeFont $g5 "tx15" "Arial"
SetShapeBox  $g5 "tx16" 9688461 4457282 397540 105768
SetShapeFont $g5 "tx16" "Arial"
SetShapeBox  $g5 "tx17" 8634912 4131109 397540 105768
SetShapeFont $g5 "tx17" "Arial"
SetShapeBox  $g5 "tx28" 3430907 2764413 77974 75989
SetShapeFont $g5 "tx28" "Arial"
SetShapeBox  $g5 "tx30" 4755306 2764413 77974 75989
SetShapeFont $g5 "tx30" "Arial"

# Slide 6 - Outcome Indicators
$s6 = $p.Slides.Item(6)
$g6 = $s6.Shapes.Item(2)
SetShapeBox  $g6 "tx17" 10903805 3248207 397540 105768
SetShapeFont $g6 "tx17" "Arial"
SetShapeBox  $g6 "tx18" 10743220 2922034 397540 105768
SetShapeFont $g6 "tx18" "Arial"
SetShapeBox  $g6 "tx19" 5598937 4544162 397540 105768
SetShapeFont $g6 "tx19" "Arial"
SetShapeBox  $g6 "tx20" 6106184 4217990 397540 105768
SetShapeFont $g6 "tx20" "Arial"
SetShapeBox  $g6 "tx21" 10633782 3896184 397540 105768
SetShapeFont $g6 "tx21" "Arial"
SetShapeBox  $g6 "tx22" 9968195 3570012 397540 105768
SetShapeFont $g6 "tx22" "Arial"
SetShapeBox  $g6 "tx38" 4012228 2434997 77974 75989
SetShapeFont $g6 "tx38" "Arial"
SetShapeBox  $g6 "tx40" 5336627 2434997 77974 75989
SetShapeFont $g6 "tx40" "Arial"
